$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prep_sheet")

# Row 19: "green olive gremolata" / "1 quart" -> "lemon caper condiment" / "3 cups"
$ws.Range("A19").Value = "lemon caper condiment"
$ws.Range("B19").Value = "3 cups"

# Row 36: clear the "2 cases" entry in column A
$ws.Range("A36").Value = ""

# Update the active sheet selection to match the saved view state
$ws.Activate()
$ws.Range("A23:B25").Select()
